$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inventario")

# Row 7: update Stock (D7) 26 -> 21
$ws.Range("D7").Value = 21

# Row 9: update Descripcion (B9), Precio (C9) and Stock (D9)
$ws.Range("B9").Value = "Fabe Naproxeno Paracetamol 10 tabletas"
$ws.Range("C9").Value = 50
$ws.Range("D9").Value = 25
